$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - new bug entry #2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "rsk"
$ws.Range("D3").Value = (Get-Date -Year 2024 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E3").Value = "Main"
$ws.Range("G3").Value = "When stopping the program, the pump continues to run."

# Row 4 - new bug entry #3 (reuses the previously-blank row 4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "rsk"
$ws.Range("D4").Value = (Get-Date -Year 2024 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E4").Value = "Main"
$ws.Range("G4").Value = "Start button pressed, fill light on; in jog mode turn fill off and back on (toggling)"

$ws.Rows.Item(4).RowHeight = 29

# Update selection to reflect the last-edited cell
$ws.Range("G4").Select()
